$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 368, shifting existing rows 368-449 down to 369-450
$ws.Rows.Item(368).Insert()

# Populate the newly inserted row 368 with the new record's data
$ws.Range("A368").Value2 = 3
$ws.Range("B368").Value2 = "Femacal de La Calera"
$ws.Range("C368").Value2 = "Coquimbo"
$ws.Range("D368").Value2 = 44889
$ws.Range("E368").Value2 = 5
$ws.Range("F368").Value2 = 100112043
$ws.Range("G368").Value2 = "Pepino ensalada"
$ws.Range("H368").Value2 = "Sin especificar"
$ws.Range("I368").Value2 = "Primera"
$ws.Range("J368").Value2 = 85
$ws.Range("K368").Value2 = 15000
$ws.Range("L368").Value2 = 16000
$ws.Range("M368").Value2 = 15471
$ws.Range("N368").Value2 = "$/caja 70 unidades"
$ws.Range("O368").Value2 = "Limache"
$ws.Range("P368").Value2 = 221
$ws.Range("Q368").Value2 = 70
$ws.Range("R368").Value2 = "Hortaliza"

# Ensure the date cell keeps the workbook's date/time number format (style index 2)
$ws.Range("D368").NumberFormat = $ws.Range("D369").NumberFormat
